$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Fitness values (column C) for rows 2 through 99
$values = @(
    13928,13376,12746,12544,12544,11212,10893,10760,10547,10190,
    9604,9471,9471,9471,9211,9211,9211,9211,9211,8899,
    8738,8738,8738,8112,8112,8112,8112,8081,8081,8081,
    8073,8073,8073,8073,8073,8073,8073,8073,8073,8073,
    8073,7708,7708,7708,7708,7708,7708,7708,7708,7708,
    7708,7708,7708,7708,7708,7708,7708,7708,7708,7708,
    7708,7708,7708,7708,7708,7708,7708,7708,7708,7590,
    7590,7590,7590,7590,7590,7590,7590,7590,7590,7590,
    7590,7590,7590,7590,7590,7590,7590,7590,7590,7590,
    7590,7590,7590,7590,7590,7590,7590,7590
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
